$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DA1 part (row 10) to the newer, more specific regulator part number,
# matching the v1.4 BOM revision ("MCP1702T" -> "MCP1702T-3302").
# A leading apostrophe is used so the cell keeps being treated as explicit text
# (preserving the existing quote-prefixed cell style) rather than Excel
# re-evaluating/clearing that formatting when the value is rewritten.
$ws.Range("B10").Value = "'MCP1702T-3302"
$ws.Range("C10").Value = "'Linear voltage regulator 3.3V"

# Widen the Description column (column C) to comfortably fit the longer text.
$ws.Columns.Item(3).ColumnWidth = 24
